$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# Row 14: fill previously-blank cells with the literal "nan" placeholder
# (matches the convention used by every other data row on this sheet).
$ws.Range("B14").Value = "nan"
$ws.Range("C14").Value = "nan"
$ws.Range("E14").Value = "nan"
$ws.Range("F14").Value = "nan"
$ws.Range("G14").Value = "nan"
$ws.Range("H14").Value = "nan"
$ws.Range("I14").Value = "nan"
$ws.Range("J14").Value = "nan"
$ws.Range("K14").Value = "nan"

# Row 15: new maintenance event added to Card2
# (A15/D15 are forced to text so "2"/"990" aren't stored as numbers, then
# the style is reset back to Normal so no extra cell format is introduced)
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2"
$ws.Range("A15").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "990"
$ws.Range("D15").Style = "Normal"

$ws.Range("L15").Value = "2\11\2025"
$ws.Range("M15").Value = "قطع سير 1270"
$ws.Range("N15").Value = "تغير سير 1270"
$ws.Range("O15").Value = "حسام"
